$wb = $excel.ActiveWorkbook

# 1. Update selection on the "validLogin" sheet (was the active/tabSelected sheet before)
$wsValid = $wb.Worksheets.Item("validLogin")
$wsValid.Range("B11").Select()

# 2. Add the new "registration" sheet at the end of the workbook
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsReg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsReg.Name = "registration"

# Header row
$wsReg.Range("A1").Value = "firstName"
$wsReg.Range("B1").Value = "middleName"
$wsReg.Range("C1").Value = "lastName"
$wsReg.Range("D1").Value = "password"
$wsReg.Range("E1").Value = "confirmPassword"

# Data row
$wsReg.Range("A2").Value = "Mayank"
$wsReg.Range("B2").Value = "Upendra"
$wsReg.Range("C2").Value = "Mishra"
$wsReg.Range("D2").Value = "User@123"
$wsReg.Range("E2").Value = "User@123"

# Hyperlinks on the password / confirmPassword cells (mirrors the other sheets' pattern)
$wsReg.Hyperlinks.Add($wsReg.Range("D2"), "mailto:User@123")
$wsReg.Hyperlinks.Add($wsReg.Range("E2"), "mailto:User@123")

# Make the new sheet's H8 the active selection (becomes the active/tabSelected sheet)
$wsReg.Range("H8").Select()
